$wb = $excel.ActiveWorkbook

# --- Step 1: on the existing "Sheet1" (AcacSuav config, the original/only
# sheet), just re-case/re-order the header row. Data rows stay the AcacSuav
# values, untouched - this sheet is kept around as the old/reference
# species config. ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "variable"
$ws1.Range("B1").Value = "value"
$ws1.Range("C1").Value = "notes"
$ws1.Range("A1:B8").Select()

# --- Step 2: add a brand-new sheet for the new PersHirs species config.
# Adding with no Before/After places it right before the active sheet, and
# Excel names it "Sheet2" - exactly matching the final tab order
# (Sheet2, Sheet1). ---
$new = $wb.Worksheets.Add()
$new.Name = "Sheet2"

# --- Step 3: fill in the new PersHirs species configuration. Header is
# lowercase variable/value, and the notes column is dropped entirely (only
# A:B used now). ---
$new.Range("A1").Value = "variable"
$new.Range("B1").Value = "value"

$new.Range("A2").Value = "species"
$new.Range("B2").Value = "PersHirs"

$new.Range("A3").Value = "dataset"
$new.Range("B3").Value = "DPers21-6117"

$new.Range("A4").Value = "raw_meta_path"
$new.Range("B4").Value = "/Users/eilishmcmaster/Documents/ReCER_base_analysis_pipeline/PersHirs/meta/samples-2023-11-07_213926/Tissue-2023-11-07_213926.csv"

$new.Range("A5").Value = "species_col_name"
$new.Range("B5").Value = "sp"

$new.Range("A6").Value = "site_col_name"
$new.Range("B6").Value = "site"

$new.Range("A7").Value = "remove_pops_less_than_n5"
$new.Range("B7").Value = $true

$new.Range("A8").Value = "downsample_pops"
$new.Range("B8").Value = $true

# match the style used on the "dataset" value cell in the original sheet
$new.Range("B3").Font.Color = $ws1.Range("B3").Font.Color

# widen column B on the new sheet to fit the long PersHirs raw_meta_path
$new.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth
$new.Columns.Item(2).ColumnWidth = 123.16666666666667

$new.Range("B11").Select()
$new.Activate()
